# "calender and CalendarTest Updated"
#
# - Rename sheet "NewEventInformation" -> "Calendar"
# - Add a new sheet "HoverOverValues" right after "Calendar", containing
#   the calendar hover/nav labels
# - Reuse existing border/fill styles where possible, and add a new
#   (slightly different) green header fill for the new sheet's header row
# - Leave the final selection on each sheet matching the saved workbook
#   state (Login!E6, Calendar!C9, HoverOverValues!E4), with
#   HoverOverValues as the active/visible tab

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsCalendar = $wb.Worksheets.Item("NewEventInformation")
$wsCalendar.Name = "Calendar"

# New third sheet, inserted immediately after "Calendar"
$wsHover = $wb.Worksheets.Add($null, $wsCalendar)
$wsHover.Name = "HoverOverValues"

# Header row
$wsHover.Range("A1").Value = "Calendar"
$wsHover.Range("B1").Value = "Companies"
$wsHover.Range("C1").Value = "Contacts"

# Nav / hover labels
$wsHover.Range("A2").Value = "New Event"
$wsHover.Range("A3").Value = "View Today"
$wsHover.Range("A4").Value = "Week View"
$wsHover.Range("A5").Value = "Month View"

# Reuse the plain bordered style already used on the Login sheet for the
# data rows (A2:C5)
$wsLogin.Range("A2").Copy() | Out-Null
$wsHover.Range("A2:C5").PasteSpecial(-4122) | Out-Null

# Reuse the bordered+filled header style from the Login sheet for the
# header row, then recolor it to the new green (FF00B050)
$wsLogin.Range("A1:C1").Copy() | Out-Null
$wsHover.Range("A1:C1").PasteSpecial(-4122) | Out-Null
$wsHover.Range("A1:C1").Interior.Color = 5287936

# Final selections / active sheet to match the saved workbook view state
$wsLogin.Range("E6").Select() | Out-Null
$wsCalendar.Range("C9").Select() | Out-Null
$wsHover.Range("E4").Select() | Out-Null

$wsHover.Activate() | Out-Null
